$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 644 (shifts old rows 644:685 -> 646:687)
$ws.Rows.Item(644).Insert()
$ws.Rows.Item(644).Insert()

# Force column A on the new rows to be stored as plain text (like the rest of
# the date column) instead of letting Excel auto-convert "yyyy/mm/dd" looking
# text into a date serial number.
$ws.Range("A644:A645").NumberFormat = "@"

# Populate the two newly inserted rows with the new daily entries
$ws.Range("A644").Value = "2026/01/17"
$ws.Range("B644").Value = "土"
$ws.Range("C644").Value = 22
$ws.Range("D644").Value = 23

$ws.Range("A645").Value = "2026/01/18"
$ws.Range("B645").Value = "日"
$ws.Range("C645").Value = 1
$ws.Range("D645").Value = 18

# Drop back to the default cell style so the new cells don't carry a
# leftover explicit "Text" number format (matches the plain, unstyled cells
# used throughout the rest of the sheet).
$ws.Range("A644:A645").Style = "Normal"

"done"
